$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = 0.75588153503657229
$ws.Range("L1").Value = 0.83808540677147003
$ws.Range("S1").Value = 0.88519648777773052
$ws.Range("AY2").Value = 0.8141520450799461
$ws.Range("BP2").Value = 0.93599627434761468
$ws.Range("Y3").Value = 0.78811613445947892
$ws.Range("AG3").Value = 0.72122078337607454
$ws.Range("E4").Value = 0.76193862847827498
$ws.Range("AF4").Value = 0.58983137582562484
$ws.Range("AJ4").Value = 0.72841411198359252
$ws.Range("C5").Value = 0.73461841691865104
$ws.Range("AA6").Value = 0.72507048242232619
$ws.Range("F7").Value = 0.78104319114547338
$ws.Range("H7").Value = 0.5484348120183512
$ws.Range("AH7").Value = 0.63016525641433985
$ws.Range("AQ7").Value = 0.74604017384975874
$ws.Range("AY8").Value = 0.68066493761630109
$ws.Range("BK8").Value = 0.94433914610430092
$ws.Range("M9").Value = 0.77377185631460543
$ws.Range("S9").Value = 0.61104941509577415
$ws.Range("E10").Value = 0.92119861813802562
$ws.Range("AM10").Value = 0.99503903978907893
$ws.Range("F11").Value = 0.62378640985514688
$ws.Range("BH11").Value = 0.68224058292860523
$ws.Range("C12").Value = 0.70855841005450826
$ws.Range("AB12").Value = 0.99081090227131197
$ws.Range("AO12").Value = 0.95055099457236469
$ws.Range("Y13").Value = 0.63392755122227573
$ws.Range("W15").Value = 0.94384584951347583
$ws.Range("AE15").Value = 0.65373891778132065
$ws.Range("BO15").Value = 0.81117057770337597
$ws.Range("N16").Value = 0.9314956651740216
$ws.Range("R16").Value = 0.83074682544909129
$ws.Range("AC16").Value = 0.62665769808669447
$ws.Range("X17").Value = 0.77027068286272293
$ws.Range("AN17").Value = 0.99771854125737369
$ws.Range("BC17").Value = 0.93106733867712865
$ws.Range("T18").Value = 0.82784647812166257
$ws.Range("AS18").Value = 0.84778974364939474
$ws.Range("AB19").Value = 0.73342381587678362
$ws.Range("AO19").Value = 0.90897268679996368
$ws.Range("K20").Value = 0.95109751129775422
$ws.Range("M20").Value = 0.76329343384408843
$ws.Range("AK20").Value = 0.60765563338858042
$ws.Range("AM21").Value = 0.81818143332677107
$ws.Range("AS21").Value = 0.94746681382366327
$ws.Range("Y23").Value = 0.74597997305783914
$ws.Range("H25").Value = 0.80410202471550929
$ws.Range("BN25").Value = 0.93213966746548138
$ws.Range("B27").Value = 0.80020662762075168
$ws.Range("Z27").Value = 0.98040820366117309
$ws.Range("I28").Value = 0.60179774493000393
$ws.Range("AC28").Value = 0.81406987504595119
$ws.Range("R29").Value = 0.85960875738654763
$ws.Range("BH30").Value = 0.81727814735126914
$ws.Range("AH31").Value = 0.94464064136158976
$ws.Range("AF33").Value = 0.60904864053247132
$ws.Range("BI33").Value = 0.79454431271775017
$ws.Range("AC34").Value = 0.98755434404007503
$ws.Range("AJ34").Value = 0.78575331416459915
$ws.Range("BB34").Value = 0.86804912280893565
$ws.Range("BJ34").Value = 0.99152206494488149
$ws.Range("AG35").Value = 0.82865266411040972
$ws.Range("AV35").Value = 0.63861201293197878
$ws.Range("AY35").Value = 0.68624356128754993
$ws.Range("AR37").Value = 0.5439011051096927
$ws.Range("J38").Value = 0.87284111847096046
$ws.Range("AM38").Value = 0.82129036156997826
$ws.Range("BA38").Value = 0.633888776148264
$ws.Range("AV39").Value = 0.86248839056447912
$ws.Range("AD40").Value = 0.67096864717974192
$ws.Range("AG41").Value = 0.89093079592349311
$ws.Range("M42").Value = 0.99854391321309133
$ws.Range("N42").Value = 0.94590389440953637
$ws.Range("AL42").Value = 0.70408440744392109
$ws.Range("AN42").Value = 0.90499081253553604
$ws.Range("Q43").Value = 0.93398985652277167
$ws.Range("AW43").Value = 0.99700562420127281
$ws.Range("BE44").Value = 0.909333013239503
$ws.Range("T45").Value = 0.76327024499487561
$ws.Range("Z45").Value = 0.87116509942907883
$ws.Range("B46").Value = 0.81861162983007385
$ws.Range("BE46").Value = 0.54480437336559429
$ws.Range("BK46").Value = 0.6626074790716614
$ws.Range("BJ47").Value = 0.98172598387759713
$ws.Range("AT48").Value = 0.81144763391836316
$ws.Range("AU48").Value = 0.63769467418178127
$ws.Range("AW48").Value = 0.71083468307237285
$ws.Range("AX48").Value = 0.58683167000737924
$ws.Range("V50").Value = 0.92413229380969342
$ws.Range("AZ50").Value = 0.93289861914123251
$ws.Range("AH51").Value = 0.64803285460248672
$ws.Range("AW51").Value = 0.55323316887286378
$ws.Range("AX51").Value = 0.68101915297976223
$ws.Range("AZ51").Value = 0.93567647871190429
$ws.Range("R52").Value = 0.97593566536988985
$ws.Range("AQ52").Value = 0.81297072199291831
$ws.Range("X53").Value = 0.75425671705206332
$ws.Range("Y53").Value = 0.74828131754028204
$ws.Range("AQ53").Value = 0.72503261517965734
$ws.Range("AZ53").Value = 0.93145491091368904
$ws.Range("N54").Value = 0.94074078302825237
$ws.Range("AE54").Value = 0.77890988564494124
$ws.Range("AZ54").Value = 0.66500364749258778
$ws.Range("BF54").Value = 0.85943223938768532
$ws.Range("E55").Value = 0.98333875709036411
$ws.Range("S55").Value = 0.76651128724249173
$ws.Range("BB55").Value = 0.67000075931924696
$ws.Range("AC56").Value = 0.92274649479067861
$ws.Range("AF56").Value = 0.89688483069189551
$ws.Range("BE56").Value = 0.94660990736368489
$ws.Range("AK57").Value = 0.9249683700452711
$ws.Range("BJ57").Value = 0.74270183309363758
$ws.Range("AX58").Value = 0.98113312051436208
$ws.Range("BM58").Value = 0.91876462399315773
$ws.Range("S59").Value = 0.68566329598108866
$ws.Range("V59").Value = 0.95051501477016132
$ws.Range("AC59").Value = 0.80528819915186878
$ws.Range("F60").Value = 0.69421907558103402
$ws.Range("I60").Value = 0.84437954563194861
$ws.Range("AL60").Value = 0.79524071530095974
$ws.Range("Q61").Value = 0.99533709884083121
$ws.Range("BK61").Value = 0.92956730318978353
$ws.Range("BK62").Value = 0.76270391188839715
$ws.Range("BP62").Value = 0.99369380654809936
$ws.Range("D63").Value = 0.88108803019031678
$ws.Range("BC63").Value = 0.82099737056204758
$ws.Range("AE64").Value = 0.78444527351596016
$ws.Range("BI64").Value = 0.91584829462785333
$ws.Range("AM65").Value = 0.9978162347119528
$ws.Range("AM66").Value = 0.57565475908514929
$ws.Range("BO66").Value = 0.85725332739055538
$ws.Range("V67").Value = 0.91928315314728948
$ws.Range("AS67").Value = 0.76455731587742859
$ws.Range("BD67").Value = 0.72868251826146269
$ws.Range("BO68").Value = 0.94347504813571981

Write-Output "Applied 136 cell updates"
